# Generate Report for Handback
#
# Row 7 of both the "zh-cn" and "de-de" sheets corresponds to the source
# file 5aa360c6-a702-40d2-98d1-6825f0d377e0.md, for which a handback was
# just received. Fill in the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns (I, J, K, P) for
# that row on each language sheet, and turn the new "Latest Target File"
# value into a hyperlink (matching the style already used for column A /
# other column-I cells on the sheet).

$wb = $excel.ActiveWorkbook

$sourceFile = "5aa360c6-a702-40d2-98d1-6825f0d377e0.md"
$latestTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97019650412f0934bddceea2b7d28088cdcf67cf/e2e/5aa360c6-a702-40d2-98d1-6825f0d377e0.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f7981e90f5bc3435f5d767fb76e607f8dd253ccc/e2e/5aa360c6-a702-40d2-98d1-6825f0d377e0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97019650412f0934bddceea2b7d28088cdcf67cf/e2e/5aa360c6-a702-40d2-98d1-6825f0d377e0.md."

function Update-HandbackRow([string]$sheetName, [string]$handbackFile, [string]$handbackDate) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Latest Target File (column I) becomes a hyperlink to the source file,
    # same as column A and the other already-linked column-I cells.
    $ws.Range("I7").Value = $sourceFile
    $ws.Hyperlinks.Add($ws.Range("I7"), $latestTargetUrl, [Type]::Missing, [Type]::Missing, $sourceFile) | Out-Null
    $ws.Range("I7").Style = $ws.Range("I5").Style

    # Latest Handback File (column J)
    $ws.Range("J7").Value = $handbackFile

    # Latest Handback DateTime (column K)
    $ws.Range("K7").Value = $handbackDate

    # Error Detail (column P)
    $ws.Range("P7").Value = $errorDetail
}

Update-HandbackRow "zh-cn" "5aa360c6-a702-40d2-98d1-6825f0d377e0.56d5c51b1841b13acd64fd9879fc4de1ce6c5eb6.zh-cn.xlf" "2016-08-25 10:57:29"
Update-HandbackRow "de-de" "5aa360c6-a702-40d2-98d1-6825f0d377e0.56d5c51b1841b13acd64fd9879fc4de1ce6c5eb6.de-de.xlf" "2016-08-25 10:57:37"
